# Apply updated invoice data (renewal terms, price-increment months,
# invoice/service start dates, price-increase percentages, and invoice
# dates) per the "Add latest code updates" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 24
$ws.Range("I2").Value = 6
$ws.Range("N2").Value = "'2025-02-15"
$ws.Range("O2").Value = "'2025-02-15"
$ws.Range("I3").Value = 7
$ws.Range("N3").Value = "'2025-04-01"
$ws.Range("O3").Value = "'2025-04-01"
$ws.Range("P3").Value = 4
$ws.Range("H4").Value = 24
$ws.Range("N4").Value = "'2025-05-10"
$ws.Range("O4").Value = "'2025-05-10"
$ws.Range("P4").Value = 5
$ws.Range("I5").Value = 6
$ws.Range("N5").Value = "'2025-06-01"
$ws.Range("O5").Value = "'2025-06-01"
$ws.Range("P5").Value = 2
$ws.Range("H6").Value = 24
$ws.Range("I6").Value = 7
$ws.Range("N6").Value = "'2025-08-01"
$ws.Range("O6").Value = "'2025-08-01"
$ws.Range("N7").Value = "'2025-09-15"
$ws.Range("O7").Value = "'2025-09-15"
$ws.Range("P7").Value = 4
$ws.Range("H8").Value = 24
$ws.Range("I8").Value = 6
$ws.Range("N8").Value = "'2025-11-01"
$ws.Range("O8").Value = "'2025-11-01"
$ws.Range("P8").Value = 5
$ws.Range("I9").Value = 7
$ws.Range("P9").Value = 2
$ws.Range("H10").Value = 24
$ws.Range("N10").Value = "'2025-02-15"
$ws.Range("O10").Value = "'2025-02-15"
$ws.Range("I11").Value = 6
$ws.Range("N11").Value = "'2025-04-01"
$ws.Range("O11").Value = "'2025-04-01"
$ws.Range("P11").Value = 4
$ws.Range("H12").Value = 24
$ws.Range("I12").Value = 7
$ws.Range("N12").Value = "'2025-05-10"
$ws.Range("O12").Value = "'2025-05-10"
$ws.Range("P12").Value = 5
$ws.Range("N13").Value = "'2025-06-01"
$ws.Range("O13").Value = "'2025-06-01"
$ws.Range("P13").Value = 2
$ws.Range("H14").Value = 24
$ws.Range("I14").Value = 6
$ws.Range("N14").Value = "'2025-08-01"
$ws.Range("O14").Value = "'2025-08-01"
$ws.Range("I15").Value = 7
$ws.Range("N15").Value = "'2025-09-15"
$ws.Range("O15").Value = "'2025-09-15"
$ws.Range("P15").Value = 4
$ws.Range("H16").Value = 24
$ws.Range("N16").Value = "'2025-11-01"
$ws.Range("O16").Value = "'2025-11-01"
$ws.Range("P16").Value = 5
$ws.Range("I17").Value = 6
$ws.Range("P17").Value = 2
$ws.Range("H18").Value = 24
$ws.Range("I18").Value = 7
$ws.Range("N18").Value = "'2025-02-15"
$ws.Range("O18").Value = "'2025-02-15"
$ws.Range("N19").Value = "'2025-04-01"
$ws.Range("O19").Value = "'2025-04-01"
$ws.Range("P19").Value = 4
$ws.Range("H20").Value = 24
$ws.Range("I20").Value = 6
$ws.Range("N20").Value = "'2025-05-10"
$ws.Range("O20").Value = "'2025-05-10"
$ws.Range("P20").Value = 5
$ws.Range("I21").Value = 7
$ws.Range("N21").Value = "'2025-06-01"
$ws.Range("O21").Value = "'2025-06-01"
$ws.Range("P21").Value = 2
$ws.Range("B34").Value = "'2025-06-15"
$ws.Range("N34").Value = "'2025-06-01"
$ws.Range("O34").Value = "'2025-06-01"
$ws.Range("B35").Value = "'2025-06-15"
$ws.Range("N35").Value = "'2025-06-01"
$ws.Range("O35").Value = "'2025-06-01"
$ws.Range("B36").Value = "'2025-06-15"
$ws.Range("N36").Value = "'2025-06-01"
$ws.Range("O36").Value = "'2025-06-01"
$ws.Range("B37").Value = "'2025-06-15"
$ws.Range("N37").Value = "'2025-06-01"
$ws.Range("O37").Value = "'2025-06-01"
$ws.Range("B38").Value = "'2025-06-15"
$ws.Range("N38").Value = "'2025-06-01"
$ws.Range("O38").Value = "'2025-06-01"
$ws.Range("B39").Value = "'2025-06-15"
$ws.Range("N39").Value = "'2025-06-01"
$ws.Range("O39").Value = "'2025-06-01"
$ws.Range("B40").Value = "'2025-06-15"
$ws.Range("N40").Value = "'2025-06-01"
$ws.Range("O40").Value = "'2025-06-01"
$ws.Range("B41").Value = "'2025-06-15"
$ws.Range("N41").Value = "'2025-06-01"
$ws.Range("O41").Value = "'2025-06-01"
$ws.Range("B42").Value = "'2025-06-15"
$ws.Range("N42").Value = "'2025-06-01"
$ws.Range("O42").Value = "'2025-06-01"
$ws.Range("B43").Value = "'2025-06-15"
$ws.Range("N43").Value = "'2025-06-01"
$ws.Range("O43").Value = "'2025-06-01"
$ws.Range("B44").Value = "'2025-06-15"
$ws.Range("N44").Value = "'2025-06-01"
$ws.Range("O44").Value = "'2025-06-01"
$ws.Range("B45").Value = "'2025-06-15"
$ws.Range("N45").Value = "'2025-06-01"
$ws.Range("O45").Value = "'2025-06-01"
$ws.Range("B46").Value = "'2025-06-15"
$ws.Range("N46").Value = "'2025-06-01"
$ws.Range("O46").Value = "'2025-06-01"
$ws.Range("B47").Value = "'2025-06-15"
$ws.Range("N47").Value = "'2025-06-01"
$ws.Range("O47").Value = "'2025-06-01"
$ws.Range("B48").Value = "'2025-06-15"
$ws.Range("N48").Value = "'2025-06-01"
$ws.Range("O48").Value = "'2025-06-01"
$ws.Range("B49").Value = "'2025-06-15"
$ws.Range("N49").Value = "'2025-06-01"
$ws.Range("O49").Value = "'2025-06-01"
$ws.Range("B50").Value = "'2025-06-15"
$ws.Range("N50").Value = "'2025-06-01"
$ws.Range("O50").Value = "'2025-06-01"
$ws.Range("B51").Value = "'2025-06-15"
$ws.Range("N51").Value = "'2025-06-01"
$ws.Range("O51").Value = "'2025-06-01"
$ws.Range("B52").Value = "'2025-06-15"
$ws.Range("N52").Value = "'2025-06-01"
$ws.Range("O52").Value = "'2025-06-01"
$ws.Range("B53").Value = "'2025-06-15"
$ws.Range("B54").Value = "'2025-06-15"
$ws.Range("B55").Value = "'2025-06-15"
$ws.Range("B56").Value = "'2025-06-15"
$ws.Range("B57").Value = "'2025-06-15"
$ws.Range("N57").Value = "'2025-06-01"
$ws.Range("O57").Value = "'2025-06-01"
